# gitbook version up
# Fill in the "filename" (E column) figure-naming scheme for rows whose
# D-column page/figure reference previously had no matching code, following
# the existing "{chapter:02d}-fig-{seq:02d}" convention already used for
# rows 3-10 (E3:E10) and rows 9-10 in column D.
#
# The order below mirrors the order the new unique strings were first
# entered (as recovered from the shared-strings table), so the shared
# string table ends up built in the same sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E19").Value = "06-fig-01"
$ws.Range("E20").Value = "07-fig-02"
$ws.Range("E24").Value = "09-fig-01"
$ws.Range("E21").Value = "-"
$ws.Range("E22").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("E25").Value = "09-fig-02"
$ws.Range("E26").Value = "09-fig-03"
$ws.Range("E27").Value = "11-fig-01"
$ws.Range("E28").Value = "11-fig-02"
$ws.Range("E30").Value = "11-fig-10"
$ws.Range("E29").Value = "11-fig-07"

# Scroll the sheet over a column and move the selection down to where the
# new entries were typed, matching the saved view state (topLeftCell="B1",
# active cell E30).
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E30").Select()

# Resize/reposition the workbook window to match the saved view.
$win = $excel.ActiveWindow
$win.Left = 8235
$win.Top = 4560
$win.Width = 12150
$win.Height = 11385
